$wb = $excel.ActiveWorkbook

# --- Summary sheet updates ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("A3").Value = 301.72000000000003
$wsSummary.Range("E3").Value = 301.72000000000003

# --- Repayment schedule sheet updates ---
$wsRepay = $wb.Worksheets.Item("Repayment schedule")

# Row 12
$wsRepay.Range("B12").Value = 14
$wsRepay.Range("C12").Value = 42157
$wsRepay.Range("F12").Value = 846.72
$wsRepay.Range("G12").Value = 1705.24
$wsRepay.Range("H12").Value = 11.75

# Row 13
$wsRepay.Range("B13").Value = 14
$wsRepay.Range("F13").Value = 850.62
$wsRepay.Range("G13").Value = 854.62
$wsRepay.Range("H13").Value = 7.85

# Row 14
$wsRepay.Range("F14").Value = 854.62
$wsRepay.Range("H14").Value = 3.93
$wsRepay.Range("K14").Value = 858.55
$wsRepay.Range("Q14").Value = 858.55

# --- Selection / view updates ---
# Summary: selection moves from D7 to D9
$wsSummary.Range("D9").Select()

# Repayment schedule: selection moves from J15 to K18, and becomes the active/selected tab
$wsRepay.Activate()
$wsRepay.Range("K18").Select()

# Input sheet: keep topLeftCell A7, but it's no longer the tab-selected sheet
# (handled implicitly since Repayment schedule becomes active)

$wb.Save()
